# Error Calculations and Plots
#
# This workbook (combination_3_ABCDF/AF/15/seed2/missing_data.xlsx) holds a
# small ID/A/B/C/D/F measurement table. This edit:
#   1. Removes two data rows entirely ("RM 232" and "SC 92"), which shifts
#      every row below them up and shrinks the used range from F35 to F33.
#   2. Updates a handful of "F" (and a couple of "B") column values on the
#      remaining rows - some previously-missing cells now carry an imputed
#      value, and a few previously-filled cells are cleared back to blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the two obsolete rows -------------------------------------
# Row 26 is "RM 232". Deleting it shifts "SC 92" up from row 28 to row 27.
$ws.Rows(26).Delete()
# "SC 92" is now at row 27; remove it too.
$ws.Rows(27).Delete()

# --- 2. Cell-level value edits on the resulting (post-delete) layout -----
$ws.Range("F19").Value = 17.81
$ws.Range("F21").Value = ""
$ws.Range("F23").Value = 16.48

$ws.Range("B26").Value = ""
$ws.Range("B27").Value = -20.4
$ws.Range("F27").Value = ""
$ws.Range("B29").Value = ""
$ws.Range("F33").Value = 17.53
